$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 37550
$ws.Range("J75").Value = 37550
$ws.Range("L75").Value = 37550
$ws.Range("N75").Value = -39422

$ws.Range("H78").Value = 37550
$ws.Range("J78").Value = 37550
$ws.Range("L78").Value = 112650
$ws.Range("N78").Value = -122010

$ws.Range("H112").Value = 566850.8
$ws.Range("J112").Value = 642321.9399999999
$ws.Range("L112").Value = 1926965.82
$ws.Range("N112").Value = -1929181.82

$ws.Range("H118").Value = 854.13336
$ws.Range("I118").Value = 696.1429000000001
$ws.Range("J118").Value = 992.375
$ws.Range("K118").Value = 2088.4287
$ws.Range("L118").Value = 2977.125
$ws.Range("M118").Value = -431.4287000000004
$ws.Range("N118").Value = -6291.125

$ws.Range("H132").Value = 247672.17
$ws.Range("I132").Value = 3646.1316
$ws.Range("J132").Value = 3338668.8
$ws.Range("K132").Value = 10938.3948
$ws.Range("L132").Value = 10016006.4
$ws.Range("M132").Value = -8408.3948
$ws.Range("N132").Value = -10021066.4

$ws.Range("H137").Value = 2503.4473
$ws.Range("I137").Value = 1437.9656
$ws.Range("K137").Value = 4313.8968
$ws.Range("M137").Value = -1763.8968

$ws.Range("H138").Value = 2089.36
$ws.Range("I138").Value = 953.2
$ws.Range("J138").Value = 2373.4
$ws.Range("K138").Value = 2859.6
$ws.Range("L138").Value = 7120.200000000001
$ws.Range("M138").Value = 2280.4
$ws.Range("N138").Value = -17400.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6701.2
$ws.Range("I74").Value = 7000.3335
$ws.Range("J74").Value = 5504.6665
$ws.Range("K74").Value = 7000.3335
$ws.Range("L74").Value = 5504.6665
$ws.Range("M74").Value = -6126.3335
$ws.Range("N74").Value = -7252.6665

$ws.Range("H77").Value = 6701.2
$ws.Range("I77").Value = 7000.3335
$ws.Range("J77").Value = 5504.6665
$ws.Range("K77").Value = 35001.6675
$ws.Range("L77").Value = 27523.3325
$ws.Range("M77").Value = -30633.6675
$ws.Range("N77").Value = -36259.3325

$ws.Range("H80").Value = 35241.5
$ws.Range("J80").Value = 35241.5
$ws.Range("L80").Value = 35241.5
$ws.Range("N80").Value = -37237.5

$ws.Range("H83").Value = 35241.5
$ws.Range("J83").Value = 35241.5
$ws.Range("L83").Value = 105724.5
$ws.Range("N83").Value = -115708.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17859920
$ws.Range("I31").Value = 1052.5
$ws.Range("J31").Value = 41671740
$ws.Range("K31").Value = 1052.5
$ws.Range("L31").Value = 41671740
$ws.Range("M31").Value = -757.5
$ws.Range("N31").Value = -41672330

$ws.Range("H34").Value = 17859920
$ws.Range("I34").Value = 1052.5
$ws.Range("J34").Value = 41671740
$ws.Range("K34").Value = 1052.5
$ws.Range("L34").Value = 41671740
$ws.Range("M34").Value = -850.5
$ws.Range("N34").Value = -41672144

$ws.Range("H64").Value = 36762.332
$ws.Range("J64").Value = 36762.332
$ws.Range("L64").Value = 36762.332
$ws.Range("N64").Value = -37258.332

$ws.Range("H67").Value = 36762.332
$ws.Range("J67").Value = 36762.332
$ws.Range("L67").Value = 36762.332
$ws.Range("N67").Value = -38478.332

$ws.Range("H115").Value = 35000
$ws.Range("J115").Value = 35000
$ws.Range("L115").Value = 35000
$ws.Range("N115").Value = -37350

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H134").Value = 12208.637
$ws.Range("I134").Value = 20884.2
$ws.Range("J134").Value = 4979
$ws.Range("K134").Value = 62652.60000000001
$ws.Range("L134").Value = 14937
$ws.Range("M134").Value = -60117.60000000001
$ws.Range("N134").Value = -20007

$ws.Range("H138").Value = 43408.89
$ws.Range("J138").Value = 43408.89
$ws.Range("L138").Value = 43408.89
$ws.Range("N138").Value = -53688.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3067.4866
$ws.Range("I122").Value = 1049.4
$ws.Range("J122").Value = 3814.926
$ws.Range("K122").Value = 9444.6
$ws.Range("L122").Value = 34334.334
$ws.Range("M122").Value = -6994.6
$ws.Range("N122").Value = -39234.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20835900
$ws.Range("I80").Value = 35716344
$ws.Range("J80").Value = 3282.4
$ws.Range("K80").Value = 35716344
$ws.Range("L80").Value = 3282.4
$ws.Range("M80").Value = -35715346
$ws.Range("N80").Value = -5278.4

$ws.Range("H83").Value = 20835900
$ws.Range("I83").Value = 35716344
$ws.Range("J83").Value = 3282.4
$ws.Range("K83").Value = 178581720
$ws.Range("L83").Value = 16412
$ws.Range("M83").Value = -178576728
$ws.Range("N83").Value = -26396

$ws.Range("H88").Value = 32693.75
$ws.Range("J88").Value = 32693.75
$ws.Range("L88").Value = 32693.75
$ws.Range("N88").Value = -33595.75

$ws.Range("H91").Value = 32693.75
$ws.Range("J91").Value = 32693.75
$ws.Range("L91").Value = 32693.75
$ws.Range("N91").Value = -35813.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H68").Value = 799.37
$ws.Range("I68").Value = 799.37
$ws.Range("K68").Value = 799.37
$ws.Range("M68").Value = -50.37

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H71").Value = 799.37
$ws.Range("I71").Value = 799.37
$ws.Range("K71").Value = 3996.85
$ws.Range("M71").Value = -252.8499999999999

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H80").Value = 47437.5
$ws.Range("J80").Value = 47437.5
$ws.Range("L80").Value = 47437.5
$ws.Range("N80").Value = -49683.5

$ws.Range("H82").Value = 1585.44
$ws.Range("J82").Value = 2204.8667
$ws.Range("L82").Value = 2204.8667
$ws.Range("N82").Value = -2926.8667

$ws.Range("H83").Value = 47437.5
$ws.Range("J83").Value = 47437.5
$ws.Range("L83").Value = 142312.5
$ws.Range("N83").Value = -153544.5

$ws.Range("H85").Value = 1585.44
$ws.Range("J85").Value = 2204.8667
$ws.Range("L85").Value = 2204.8667
$ws.Range("N85").Value = -4700.8667

$ws.Range("H88").Value = 36000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 36000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 36000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -36856

$ws.Range("H91").Value = 36000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 36000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 36000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -38964

$ws.Range("H122").Value = 3651.8462
$ws.Range("I122").Value = 2037.8
$ws.Range("J122").Value = 9032
$ws.Range("K122").Value = 6113.4
$ws.Range("L122").Value = 27096
$ws.Range("M122").Value = -3663.4
$ws.Range("N122").Value = -31996

$ws.Range("H136").Value = 3688.9412
$ws.Range("J136").Value = 6038
$ws.Range("L136").Value = 18114
$ws.Range("N136").Value = -23214

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 18700
$ws.Range("J57").Value = 18700
$ws.Range("L57").Value = 18700
$ws.Range("N57").Value = -20208

$ws.Range("H132").Value = 13335549
$ws.Range("I132").Value = 1335.6875
$ws.Range("J132").Value = 37040816
$ws.Range("K132").Value = 4007.0625
$ws.Range("L132").Value = 111122448
$ws.Range("M132").Value = -1477.0625
$ws.Range("N132").Value = -111127508
